# Apply the renamed worksheet tabs and make the first sheet ("Reference
# Guide") the active tab, matching the committed XML diff:
#   - bookViews/workbookView gains activeTab="0"
#   - Flows_NEWTSF      -> Flows_New TSF
#   - Flows_OLDTSF      -> Flows_Old TSF
#   - Flows_UG2P        -> Flows_UG2 Plant
#   - Flows_UG2S        -> Flows_UG2 Main
#   - Flows_UG2N        -> Flows_UG2 North
#   - Flows_MERP        -> Flows_Merensky Plant
#   - Flows_MERS        -> Flows_Merensky South
#   - Flows_STOCKPILE   -> Flows_Stockpile1

$wb = $excel.ActiveWorkbook

$renames = @{
    "Flows_NEWTSF"    = "Flows_New TSF"
    "Flows_OLDTSF"    = "Flows_Old TSF"
    "Flows_UG2P"      = "Flows_UG2 Plant"
    "Flows_UG2S"      = "Flows_UG2 Main"
    "Flows_UG2N"      = "Flows_UG2 North"
    "Flows_MERP"      = "Flows_Merensky Plant"
    "Flows_MERS"      = "Flows_Merensky South"
    "Flows_STOCKPILE" = "Flows_Stockpile1"
}

foreach ($ws in $wb.Worksheets) {
    $oldName = $ws.Name
    if ($renames.ContainsKey($oldName)) {
        $ws.Name = $renames[$oldName]
    }
}

# Make the first sheet the active tab (workbookView activeTab="0").
$wb.Worksheets.Item(1).Activate()
